$d = $word.ActiveDocument

$replacements = @(
    @("405÷4=", "458÷3="),
    @("353÷4=", "141÷3="),
    @("772÷7=", "137÷2="),
    @("457÷5=", "329÷6="),
    @("357÷4=", "764÷2="),
    @("268÷9=", "750÷2="),
    @("243÷6=", "611÷5="),
    @("851÷8=", "845÷6="),
    @("249÷6=", "722÷5="),
    @("605÷7=", "222÷2="),
    @("765÷2=", "995÷6="),
    @("635÷6=", "999÷5="),
    @("876÷3=", "370÷4="),
    @("871÷6=", "884÷9="),
    @("802÷7=", "236÷5="),
    @("565÷2=", "782÷7="),
    @("818÷5=", "524÷8="),
    @("773÷5=", "809÷5="),
    @("694÷7=", "276÷4="),
    @("191÷9=", "514÷2="),
    @("972÷3=", "135÷5="),
    @("940÷5=", "314÷5="),
    @("341÷6=", "451÷4="),
    @("979÷7=", "608÷7="),
    @("879÷8=", "319÷9=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
